# Fix rows that had been polluted by other tickers' rows mixed into this
# sheet: recompute open_price/close_price/high_price/low_price/
# shares_outstanding (cols D:H) for every historical row, and renormalize
# fixed_ticker (col I) back to this company's ticker (MRVL) on every row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 43,6

$data[0,0] = 13.92315180070411; $data[0,1] = 14.50815868377686; $data[0,2] = 15.10216501249978; $data[0,3] = 13.41014563273013; $data[0,4] = 862100000; $data[0,5] = "MRVL"
$data[1,0] = 12.7364545970468; $data[1,1] = 12.63709163665772; $data[1,2] = 13.1339038542565; $data[1,3] = 11.78799506553074; $data[1,4] = 862100000; $data[1,5] = "MRVL"
$data[2,0] = 11.25852150735304; $data[2,1] = 10.22429847717285; $data[2,2] = 12.11130284173944; $data[2,3] = 9.072136700872576; $data[2,4] = 862100000; $data[2,5] = "MRVL"
$data[3,0] = 7.497491551601611; $data[3,1] = 8.091080665588379; $data[3,2] = 8.237195726346128; $data[3,3] = 7.479226842415916; $data[3,4] = 862100000; $data[3,5] = "MRVL"
$data[4,0] = 8.01251383965249; $data[4,1] = 8.785247802734375; $data[4,2] = 8.987630716233259; $data[4,3] = 7.653744248991547; $data[4,4] = 862100000; $data[4,5] = "MRVL"
$data[5,0] = 9.17914328880609; $data[5,1] = 9.465991020202637; $data[5,2] = 9.521510419881508; $data[5,3] = 8.642459582531012; $data[5,4] = 862100000; $data[5,5] = "MRVL"
$data[6,0] = 10.98296741527677; $data[6,1] = 11.5414228439331; $data[6,2] = 11.72757554112706; $data[6,3] = 10.48966489805321; $data[6,4] = 862100000; $data[6,5] = "MRVL"
$data[7,0] = 12.32425280953912; $data[7,1] = 13.40893650054932; $data[7,2] = 14.02608405634792; $data[7,3] = 11.50138910455618; $data[7,4] = 862100000; $data[7,5] = "MRVL"
$data[8,0] = 14.5466996624115; $data[8,1] = 14.65000152587891; $data[8,2] = 15.30737229842711; $data[8,3] = 14.20862303480435; $data[8,4] = 862100000; $data[8,5] = "MRVL"
$data[9,0] = 14.16917142126092; $data[9,1] = 16.25259590148926; $data[9,2] = 16.78052195403242; $data[9,3] = 13.80150885943175; $data[9,4] = 862100000; $data[9,5] = "MRVL"
$data[10,0] = 14.74968568316743; $data[10,1] = 16.94463539123535; $data[10,2] = 17.00140083415686; $data[10,3] = 14.06849405221619; $data[10,4] = 862100000; $data[10,5] = "MRVL"
$data[11,0] = 17.62663225560296; $data[11,1] = 21.20511436462402; $data[11,2] = 22.98960840836444; $data[11,3] = 16.94320935149147; $data[11,4] = 862100000; $data[11,5] = "MRVL"
$data[12,0] = 21.99596029275611; $data[12,1] = 22.35764122009277; $data[12,2] = 22.79546741150349; $data[12,3] = 19.38804490651938; $data[12,4] = 862100000; $data[12,5] = "MRVL"
$data[13,0] = 18.97489053677973; $data[13,1] = 20.5593147277832; $data[13,2] = 21.71422424220528; $data[13,3] = 18.54537950121776; $data[13,4] = 862100000; $data[13,5] = "MRVL"
$data[14,0] = 20.38756248064068; $data[14,1] = 19.79412269592285; $data[14,2] = 21.29686767335567; $data[14,3] = 18.64552685519263; $data[14,4] = 862100000; $data[14,5] = "MRVL"
$data[15,0] = 15.78482322052228; $data[15,1] = 15.46797466278076; $data[15,2] = 16.95620234855605; $data[15,3] = 14.10456441321453; $data[15,4] = 862100000; $data[15,5] = "MRVL"
$data[16,0] = 17.88279809726178; $data[16,1] = 19.2324447631836; $data[16,2] = 19.5505753427233; $data[16,3] = 17.30438019900461; $data[16,4] = 862100000; $data[16,5] = "MRVL"
$data[17,0] = 24.32712140668389; $data[17,1] = 21.56179618835449; $data[17,2] = 24.49149392338618; $data[17,3] = 20.59489946521377; $data[17,4] = 862100000; $data[17,5] = "MRVL"
$data[18,0] = 25.54847386198296; $data[18,1] = 23.23205184936523; $data[18,2] = 25.96523605102153; $data[18,3] = 22.46637398350357; $data[18,4] = 862100000; $data[18,5] = "MRVL"
$data[19,0] = 23.85657126427364; $data[19,1] = 25.62516593933105; $data[19,2] = 27.57839077391026; $data[19,3] = 23.69137293184697; $data[19,4] = 862100000; $data[19,5] = "MRVL"
$data[20,0] = 23.67695701689456; $data[20,1] = 20.74533843994141; $data[20,2] = 25.37164647371204; $data[20,3] = 18.99221192811104; $data[20,4] = 862100000; $data[20,5] = "MRVL"
$data[21,0] = 25.45195423599868; $data[21,1] = 31.8588924407959; $data[21,2] = 32.19095953600858; $data[21,3] = 24.39715342263057; $data[21,4] = 862100000; $data[21,5] = "MRVL"
$data[22,0] = 36.20344458897831; $data[22,1] = 37.93487167358398; $data[22,2] = 38.39463022322374; $data[22,3] = 31.82107705205154; $data[22,4] = 862100000; $data[22,5] = "MRVL"
$data[23,0] = 37.1263192134339; $data[23,1] = 45.34504699707031; $data[23,2] = 45.43320990629748; $data[23,3] = 35.7646912294087; $data[23,4] = 862100000; $data[23,5] = "MRVL"
$data[24,0] = 51.06209958854041; $data[24,1] = 47.35455322265625; $data[24,2] = 53.51417883563862; $data[24,3] = 45.25557398579235; $data[24,4] = 862100000; $data[24,5] = "MRVL"
$data[25,0] = 44.66143494604876; $data[25,1] = 47.43067932128906; $data[25,2] = 48.32430068002208; $data[25,3] = 40.05584851966153; $data[25,4] = 862100000; $data[25,5] = "MRVL"
$data[26,0] = 59.96558231551443; $data[26,1] = 60.15235900878906; $data[26,2] = 62.98352197795027; $data[26,3] = 56.13171571665011; $data[26,4] = 862100000; $data[26,5] = "MRVL"
$data[27,0] = 67.56910882441639; $data[27,1] = 70.02902984619141; $data[27,2] = 74.89967797355587; $data[27,3] = 66.68353575516073; $data[27,4] = 862100000; $data[27,5] = "MRVL"
$data[28,0] = 71.12276401627045; $data[28,1] = 67.28254699707031; $data[28,2] = 76.99140016082971; $data[28,3] = 60.02552163665128; $data[28,4] = 862100000; $data[28,5] = "MRVL"
$data[29,0] = 57.27088475095857; $data[29,1] = 58.29586791992188; $data[29,2] = 62.19868410988307; $data[29,3] = 49.68207298629233; $data[29,4] = 862100000; $data[29,5] = "MRVL"
$data[30,0] = 54.26037938518176; $data[30,1] = 46.20718383789063; $data[30,2] = 57.66522390746032; $data[30,3] = 45.7433343515468; $data[30,4] = 862100000; $data[30,5] = "MRVL"
$data[31,0] = 40.33576061204968; $data[31,1] = 45.96814346313477; $data[31,2] = 46.05707596281563; $data[31,3] = 36.2745160756971; $data[31,4] = 862100000; $data[31,5] = "MRVL"
$data[32,0] = 42.87637951910557; $data[32,1] = 44.68764114379883; $data[32,2] = 48.98319735262233; $data[32,3] = 41.7084630538748; $data[32,4] = 862100000; $data[32,5] = "MRVL"
$data[33,0] = 39.33194894620095; $data[33,1] = 57.97695922851562; $data[33,2] = 67.39362661953609; $data[33,3] = 38.29115973077777; $data[33,4] = 862100000; $data[33,5] = "MRVL"
$data[34,0] = 64.36557211446957; $data[34,1] = 57.79704666137695; $data[34,2] = 65.77452834502068; $data[34,3] = 51.84370279926088; $data[34,4] = 862100000; $data[34,5] = "MRVL"
$data[35,0] = 46.73673911681581; $data[35,1] = 55.35894775390625; $data[35,2] = 57.73303446537559; $data[35,3] = 46.02153547229005; $data[35,4] = 862100000; $data[35,5] = "MRVL"
$data[36,0] = 67.81761713948933; $data[36,1] = 71.25822448730469; $data[36,2] = 72.45149336187964; $data[36,3] = 63.17380244688746; $data[36,4] = 862100000; $data[36,5] = "MRVL"
$data[37,0] = 64.18231575750241; $data[37,1] = 68.48170471191406; $data[37,2] = 78.06576479268197; $data[37,3] = 63.03780469122229; $data[37,4] = 862100000; $data[37,5] = "MRVL"
$data[38,0] = 65.29154673026865; $data[38,1] = 75.93939208984375; $data[38,2] = 78.68851184672683; $data[38,3] = 52.9802764185465; $data[38,4] = 862100000; $data[38,5] = "MRVL"
$data[39,0] = 83.93743212148811; $data[39,1] = 92.40096282958984; $data[39,2] = 95.25204496497716; $data[39,3] = 82.85083477065204; $data[39,4] = 862100000; $data[39,5] = "MRVL"
$data[40,0] = 108.4559464099061; $data[40,1] = 91.58014678955078; $data[40,2] = 121.4918045515922; $data[40,3] = 85.90500781602866; $data[40,4] = 862100000; $data[40,5] = "MRVL"
$data[41,0] = 60.061037093219; $data[41,1] = 60.10097503662109; $data[41,2] = 67.49002783373193; $data[41,3] = 53.7004558539689; $data[41,4] = 862100000; $data[41,5] = "MRVL"
$data[42,0] = 77.01900468823257; $data[42,1] = 62.82839965820312; $data[42,2] = 80.00702415624012; $data[42,3] = 62.57856507602716; $data[42,4] = 862100000; $data[42,5] = "MRVL"

$ws.Range("D2:I44").Value = $data
